# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 ("R") updated totals
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 411
$wsOff.Range("C3").Value = 308
$wsOff.Range("D3").Value = 97
$wsOff.Range("E3").Value = 45

# DEF sheet - row 3 ("R") updated totals
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 408
$wsDef.Range("C3").Value = 268
$wsDef.Range("D3").Value = 103
$wsDef.Range("E3").Value = 54
$wsDef.Range("F3").Value = 7
